$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

$ws.Range("D2").Value = "26.711.70"
$ws.Range("E2").Value = "  -1.12%  "
$ws.Range("D3").Value = "1.599.33"
$ws.Range("E3").Value = "  -1.35%  "
$ws.Range("E4").Value = "  -0.06%  "
Set-TextValue $ws.Range("D5") "211.87"
$ws.Range("E5").Value = "  -0.79%  "
Set-TextValue $ws.Range("D6") "0.513"
$ws.Range("E6").Value = "  +0.51%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  -1.33%  "
$ws.Range("E9").Value = "  -1.47%  "
Set-TextValue $ws.Range("D10") "19.75"
$ws.Range("E10").Value = "  -0.84%  "
$ws.Range("E11").Value = "  -0.05%  "
$ws.Range("D12").Value = "1.825.76"
$ws.Range("E12").Value = "  -1.20%  "
$ws.Range("D13").Value = "1.597.31"
$ws.Range("E13").Value = "  -1.14%  "
$ws.Range("E14").Value = "  -1.99%  "
$ws.Range("E15").Value = "  -2.27%  "
Set-TextValue $ws.Range("D16") "65.09"
$ws.Range("E16").Value = "  +1.28%  "
$ws.Range("D17").Value = "26.694.27"
$ws.Range("E17").Value = "  -1.18%  "
$ws.Range("D18").Value = "0.0₃0729"
$ws.Range("E18").Value = "  -0.85%  "
Set-TextValue $ws.Range("D19") "210.21"
$ws.Range("E19").Value = "  -1.36%  "
$ws.Range("E20").Value = "  -0.03%  "
Set-TextValue $ws.Range("D21") "6.71"
$ws.Range("E21").Value = "  -1.32%  "
$ws.Range("E22").Value = "  -1.49%  "
Set-TextValue $ws.Range("D23") "2.30"
$ws.Range("E23").Value = "  -2.23%  "
$ws.Range("E24").Value = "  -0.32%  "
Set-TextValue $ws.Range("D25") "146.98"
$ws.Range("E25").Value = "  -0.17%  "
$ws.Range("E26").Value = "  -0.21%  "
$ws.Range("E27").Value = "  -4.17%  "
Set-TextValue $ws.Range("D28") "0.115"
$ws.Range("E28").Value = "  +1.11%  "
Set-TextValue $ws.Range("D29") "15.32"
$ws.Range("E29").Value = "  -1.01%  "
$ws.Range("E30").Value = "  -1.04%  "
$ws.Range("E31").Value = "  -1.10%  "
$ws.Range("E32").Value = "  -1.78%  "
Set-TextValue $ws.Range("D33") "0.670"
$ws.Range("E33").Value = "  -5.29%  "
Set-TextValue $ws.Range("D34") "2.91"
$ws.Range("D35").Value = "1.296.63"
$ws.Range("E35").Value = "  -2.81%  "
Set-TextValue $ws.Range("D36") "2.44"
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("E37").Value = "  -5.17%  "
$ws.Range("E38").Value = "  -2.17%  "
$ws.Range("E39").Value = "  +0.65%  "
$ws.Range("E40").Value = "  +0.00%  "
Set-TextValue $ws.Range("D41") "0.791"
$ws.Range("E41").Value = "  -0.67%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D42") "5.37"
$ws.Range("E42").Value = "  +0.68%  "
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D43") "2.20"
$ws.Range("E43").Value = "  -0.22%  "
Set-TextValue $ws.Range("D44") "64.00"
$ws.Range("E44").Value = "  +0.30%  "
$ws.Range("D45").Value = "1.737.15"
$ws.Range("E45").Value = "  -1.26%  "
Set-TextValue $ws.Range("D46") "0.885"
$ws.Range("E46").Value = "  +4.31%  "
Set-TextValue $ws.Range("D47") "89.97"
$ws.Range("E47").Value = "  +0.19%  "
$ws.Range("E48").Value = "  +0.29%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D49") "0.0986"
$ws.Range("E49").Value = "  -0.51%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D50") "0.0504"
$ws.Range("E50").Value = "  -2.23%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D51") "7.49"
$ws.Range("E51").Value = "  -1.20%  "
